{"js": "// Lattice-multiplication exercise table: replace each cell's problem with\n// a new one, keeping the table's row/column layout and run formatting\n// (sz=32) intact. Each cell's text is five logical lines joined by the\n// vertical-tab line-break character (\\u000b), matching how Word exposes\n// <w:br/> separated runs through body.values/getRange text.\nconst VT = \"\\u000b\";\n\n// Cells are addressed row-major (row, col), 0-based, matching the\n// pre-edit layout of 5 rows x 3 columns.\nconst replacements = [\n  // row 0\n  { row: 0, col: 0, top: \"62 x 94\", mid: \"  9    4\", l1: \"6|    |\", l2: \"2|    |\" },\n  { row: 0, col: 1, top: \"68 x 29\", mid: \"  2    9\", l1: \"6|    |\", l2: \"8|    |\" },\n  { row: 0, col: 2, top: \"49 x 22\", mid: \"  2    2\", l1: \"4|    |\", l2: \"9|    |\" },\n  // row 1\n  { row: 1, col: 0, top: \"17 x 53\", mid: \"  5    3\", l1: \"1|    |\", l2: \"7|    |\" },\n  { row: 1, col: 1, top: \"40 x 23\", mid: \"  2    3\", l1: \"4|    |\", l2: \"0|    |\" },\n  { row: 1, col: 2, top: \"97 x 64\", mid: \"  6    4\", l1: \"9|    |\", l2: \"7|    |\" },\n  // row 2\n  { row: 2, col: 0, top: \"41 x 22\", mid: \"  2    2\", l1: \"4|    |\", l2: \"1|    |\" },\n  { row: 2, col: 1, top: \"30 x 86\", mid: \"  8    6\", l1: \"3|    |\", l2: \"0|    |\" },\n  { row: 2, col: 2, top: \"24 x 13\", mid: \"  1    3\", l1: \"2|    |\", l2: \"4|    |\" },\n  // row 3\n  { row: 3, col: 0, top: \"42 x 36\", mid: \"  3    6\", l1: \"4|    |\", l2: \"2|    |\" },\n  { row: 3, col: 1, top: \"59 x 50\", mid: \"  5    0\", l1: \"5|    |\", l2: \"9|    |\" },\n  { row: 3, col: 2, top: \"78 x 51\", mid: \"  5    1\", l1: \"7|    |\", l2: \"8|    |\" },\n  // row 4\n  { row: 4, col: 0, top: \"91 x 98\", mid: \"  9    8\", l1: \"9|    |\", l2: \"1|    |\" },\n  { row: 4, col: 1, top: \"82 x 62\", mid: \"  6    2\", l1: \"8|    |\", l2: \"2|    |\" },\n  { row: 4, col: 2, top: \"27 x 93\", mid: \"  9    3\", l1: \"2|    |\", l2: \"7|    |\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const text = [r.top, r.mid, \"  ----\", r.l1, r.l2].join(VT);\n  const range = cell.body.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication exercise table: replace each cell's problem with\n# a new one, keeping the table's row/column layout and run formatting\n# (sz=32) intact. [char]11 is the vertical-tab line-break character that\n# Word COM uses for <w:br/> inside Range.Text.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$vt = [char]11\n\n$replacements = @(\n    # row, col, top line, middle line, bottom-left digit line, bottom-right digit line\n    @{ Row = 1; Col = 1; Top = \"62 x 94\"; Mid = \"  9    4\"; L1 = \"6|    |\"; L2 = \"2|    |\" },\n    @{ Row = 1; Col = 2; Top = \"68 x 29\"; Mid = \"  2    9\"; L1 = \"6|    |\"; L2 = \"8|    |\" },\n    @{ Row = 1; Col = 3; Top = \"49 x 22\"; Mid = \"  2    2\"; L1 = \"4|    |\"; L2 = \"9|    |\" },\n    @{ Row = 2; Col = 1; Top = \"17 x 53\"; Mid = \"  5    3\"; L1 = \"1|    |\"; L2 = \"7|    |\" },\n    @{ Row = 2; Col = 2; Top = \"40 x 23\"; Mid = \"  2    3\"; L1 = \"4|    |\"; L2 = \"0|    |\" },\n    @{ Row = 2; Col = 3; Top = \"97 x 64\"; Mid = \"  6    4\"; L1 = \"9|    |\"; L2 = \"7|    |\" },\n    @{ Row = 3; Col = 1; Top = \"41 x 22\"; Mid = \"  2    2\"; L1 = \"4|    |\"; L2 = \"1|    |\" },\n    @{ Row = 3; Col = 2; Top = \"30 x 86\"; Mid = \"  8    6\"; L1 = \"3|    |\"; L2 = \"0|    |\" },\n    @{ Row = 3; Col = 3; Top = \"24 x 13\"; Mid = \"  1    3\"; L1 = \"2|    |\"; L2 = \"4|    |\" },\n    @{ Row = 4; Col = 1; Top = \"42 x 36\"; Mid = \"  3    6\"; L1 = \"4|    |\"; L2 = \"2|    |\" },\n    @{ Row = 4; Col = 2; Top = \"59 x 50\"; Mid = \"  5    0\"; L1 = \"5|    |\"; L2 = \"9|    |\" },\n    @{ Row = 4; Col = 3; Top = \"78 x 51\"; Mid = \"  5    1\"; L1 = \"7|    |\"; L2 = \"8|    |\" },\n    @{ Row = 5; Col = 1; Top = \"91 x 98\"; Mid = \"  9    8\"; L1 = \"9|    |\"; L2 = \"1|    |\" },\n    @{ Row = 5; Col = 2; Top = \"82 x 62\"; Mid = \"  6    2\"; L1 = \"8|    |\"; L2 = \"2|    |\" },\n    @{ Row = 5; Col = 3; Top = \"27 x 93\"; Mid = \"  9    3\"; L1 = \"2|    |\"; L2 = \"7|    |\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $cell.Range.Text = $r.Top + $vt + $r.Mid + $vt + \"  ----\" + $vt + $r.L1 + $vt + $r.L2\n}\n"}
